$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with new column names
$ws.Range("A1").Value = "trade_date"
$ws.Range("B1").Value = "value_date"
$ws.Range("C1").Value = "currency"
$ws.Range("D1").Value = "underlying_currency"
$ws.Range("E1").Value = "counter_currency"
$ws.Range("F1").Value = "spot_price"
$ws.Range("G1").Value = "fixing_level"
$ws.Range("H1").Value = "trade_id"

# Copy the formatting/style from the old C1 header cell (bold, border, centered)
# to the newly added header cells D1:H1 so they match A1:C1 formatting
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update data row (row 2) with the extracted trade field values
$ws.Range("A2").Value = "28/09/2023"
$ws.Range("B2").Value = "30/09/2023"
$ws.Range("C2").Value = "USD"
$ws.Range("D2").Value = "EUR"
$ws.Range("E2").Value = "USD"

# Numeric-looking values must stay as literal text (not be coerced to
# floating point numbers, which would introduce rounding artifacts), so
# force them in as text via a leading quote prefix, then reset the style
# back to Normal so no extra numeric/text style gets attached to the cell.
$ws.Range("F2").Value = "'1.0523"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'1.0535"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "FX20230928001"
